{"js": "// Clarify the Core-h / GPU-h allocation labels by stating which node\n// type each refers to, matching commit \"Clarifies coreh and GPUh\n// definitions\":\n//   \"Total Core-h:\"  -> \"Total Core-h for CPU node use:\"\n//   \"Total GPU-h: \"  -> \"Total GPU-h for GPU node use: \"\n\nconst body = context.document.body;\n\n// --- \"Total Core-h\" -> \"Total Core-h for CPU node use\" ---\nconst coreResults = body.search(\"Total Core-h\", { matchCase: true, matchWholeWord: false });\ncoreResults.load(\"text\");\nawait context.sync();\n\nif (coreResults.items.length > 0) {\n  coreResults.items[0].insertText(\" for CPU node use\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// --- \"Total GPU-h\" -> \"Total GPU-h for GPU node use\" ---\nconst gpuResults = body.search(\"Total GPU-h\", { matchCase: true, matchWholeWord: false });\ngpuResults.load(\"text\");\nawait context.sync();\n\nif (gpuResults.items.length > 0) {\n  gpuResults.items[0].insertText(\" for GPU node use\", Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# Clarify the Core-h / GPU-h allocation labels by stating which node\n# type each refers to, matching commit \"Clarifies coreh and GPUh\n# definitions\":\n#   \"Total Core-h:\"  -> \"Total Core-h for CPU node use:\"\n#   \"Total GPU-h: \"  -> \"Total GPU-h for GPU node use: \"\n\n$d = $word.ActiveDocument\n\n# --- \"Total Core-h\" -> \"Total Core-h for CPU node use\" ---\n$findCpu = $d.Content.Find\n$findCpu.ClearFormatting()\n$findCpu.MatchCase = $true\n$findCpu.MatchWholeWord = $false\n$findCpu.MatchWildcards = $false\n$findCpu.Text = \"Total Core-h\"\n$findCpu.Replacement.ClearFormatting()\n$findCpu.Replacement.Text = \"Total Core-h for CPU node use\"\n$findCpu.Execute(\n    $findCpu.Text, $true, $false, $false, $false, $false, $true, 1, $false,\n    $findCpu.Replacement.Text, 2\n) | Out-Null\n\n# --- \"Total GPU-h\" -> \"Total GPU-h for GPU node use\" ---\n$findGpu = $d.Content.Find\n$findGpu.ClearFormatting()\n$findGpu.MatchCase = $true\n$findGpu.MatchWholeWord = $false\n$findGpu.MatchWildcards = $false\n$findGpu.Text = \"Total GPU-h\"\n$findGpu.Replacement.ClearFormatting()\n$findGpu.Replacement.Text = \"Total GPU-h for GPU node use\"\n$findGpu.Execute(\n    $findGpu.Text, $true, $false, $false, $false, $false, $true, 1, $false,\n    $findGpu.Replacement.Text, 2\n) | Out-Null\n"}
